# papers/pldi2011/drawings.pptx, slide 4 ("section 4" figures):
#   - the byte-range table's row header "X" becomes lower-case "x"
#   - the owner-id table's row header "Y " becomes lower-case "y" (and the
#     trailing space is dropped)
#   - that same table's second header line "changes" is capitalized to
#     "Changes"
# Commit message: "adding setVariableChangeHook to section 4."

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(4)

# --- Table 4 (the "X" / changes row-header cell) ---------------------------
$xTable = $slide.Shapes.Item(3).Table
$xCell = $xTable.Cell(1, 1)
$xCell.Shape.TextFrame.TextRange.Text = "x`rchanges"

# --- Table 6 (the "Y " / changes row-header cell) ---------------------------
$yTable = $slide.Shapes.Item(7).Table
$yCell = $yTable.Cell(1, 1)
$yCell.Shape.TextFrame.TextRange.Text = "y`rChanges"
